# Developpement de solve et bugfix
# Update capacity values on the "Camions" sheet and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Camions")

# Bugfix: update truck capacity values (column C)
$ws.Range("C2").Value = 150
$ws.Range("C3").Value = 100
$ws.Range("C4").Value = 85

# Leave the selection where the user last left it when saving
$ws.Range("C4").Select() | Out-Null
